# Fix typo and naming description in columns of template sheets
#
# - "QBIC sample ids*" / "QBIC sample ids" -> singular "QBIC sample id*" / "QBIC sample id"
# - Shorten the overly long mandatory-field description
# - Remove the "single-end"/"paired-end" example values from the
#   "Allowed-Values" sheet (the Sequencing read type column no longer
#   shows those literal examples) and restyle the header/example rows
#   (top/bottom accent borders) on that sheet.

$wb = $excel.ActiveWorkbook

$meta    = $wb.Worksheets.Item("Metadata")
$prop    = $wb.Worksheets.Item("Property information")
$allowed = $wb.Worksheets.Item("Allowed-Values")

# ---------------------------------------------------------------------
# 1. Text fixes
# ---------------------------------------------------------------------

# "Metadata" sheet header
$meta.Range("A1").Value = "QBIC sample id*"

# "Property information" sheet
$prop.Range("A2").Value = "QBIC sample id"
$prop.Range("C2").Value = "Each measurement need to be linked to at least on analyte sample."

# "Allowed-Values" sheet
$allowed.Range("A1").Value = "QBIC sample id*"
$allowed.Range("F2").Value = "Free text"
$allowed.Range("F3").ClearContents()

# ---------------------------------------------------------------------
# 2. Restyle the "Allowed-Values" header/example rows
#    (accent border colour switches from blue (index 10) to
#    green (index 11) on the inner edge between row 1 / row 2)
# ---------------------------------------------------------------------

$xlEdgeTop = 8
$xlEdgeBottom = 9
$accentGreen = 65280

$row1Cols = @("A","B","C","D","E","G","H","I","J","K","L","M")
foreach ($col in $row1Cols) {
    $allowed.Range($col + "1").Borders.Item($xlEdgeBottom).Color = $accentGreen
}

$row2Cols = @("A","B","C","D","E","G","H","I","J","K","L","M")
foreach ($col in $row2Cols) {
    $allowed.Range($col + "2").Borders.Item($xlEdgeTop).Color = $accentGreen
}
